# Generate Report for Handoff
#
# The "b.md" source file has now been handed off for localization, so its
# row on each sheet moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and the zh-cn / de-de sheets record the new handoff
# package (file name + timestamp) that was generated for it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: b.md row (row 3) status for both locales -------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: b.md row (row 3) -----------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-09 15:26:04"

# Keep the hyperlink's visible text in sync with the new handoff file name.
$zhCnLink = $wsZhCn.Range("C3").Hyperlinks.Item(1)
$zhCnLink.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# --- de-de sheet: b.md row (row 3) -----------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-09 15:26:12"

$deDeLink = $wsDeDe.Range("C3").Hyperlinks.Item(1)
$deDeLink.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
